# "Chore: all excercises complete"
# Fill in the "Quantity Check" column (G) with ISBLANK() formulas that
# flag rows where the Quantity (column E) is empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# G2:G30 -> =ISBLANK(E<row>)
for ($row = 2; $row -le 30; $row++) {
    $ws.Range("G$row").Formula = "=ISBLANK(E$row)"
}

# Reflect the reviewer's final on-screen state: zoomed in a bit more and
# the active cell resting on the newly-completed G2 formula.
$ws.Activate()
$excel.ActiveWindow.Zoom = 145
$ws.Range("G2").Select()
